$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.944.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.31%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.648.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.13%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.20%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'216.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.31%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.5068"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.04%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.28%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.61%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06456"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.96%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'20.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +6.10%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07819"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.79%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.99%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'1.873.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.01%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'WrappedEther"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'1.645.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.74%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.5630"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.30%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0₅7738"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.13%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'63.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.34%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'25.955.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.28%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.33%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'194.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.15%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.380"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.58%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'9.973"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.99%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.144"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.96%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.23%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.804"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -5.75%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'141.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.47%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.1237"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.92%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'6.833"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.12%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'15.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.23%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.249"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.90%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.04974"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.78%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.312"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.94%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.248"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.79%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +2.48%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.390"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.89%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.9091"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.53%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.5583"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.04%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.565"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.96%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.132.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.31%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.98%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.37%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'5.530"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.02%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.8032"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.62%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'98.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.69%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.783.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.04%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0₈111"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -7.09%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'55.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.19%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.4288"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'7.769"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.46%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.05053"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.59%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.9989"
$ws.Range("D51").Style = "Normal"
